$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1-15")

# Narrow columns AD (30) and AE (31) to match the other holiday/weekend columns
$ws.Columns.Item(30).ColumnWidth = 1.67
$ws.Columns.Item(31).ColumnWidth = 1.67

# Copy the "SAT" column (Z:AA) formatting onto AD:AE for rows 2-27,
# since Feb 15 2021 (the Monday in columns AD/AE) is now a federal holiday
# (Presidents' Day) and should be highlighted the same way weekend columns are.
$src = $ws.Range("Z2:AA27")
$dst = $ws.Range("AD2:AE27")
$src.Copy()
$dst.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Mark attendance "X" in the site rows, matching the weekend columns
$xRows = @(5,6,8,9,11,12,14,15,17,18,20,21,23,24,26,27)
foreach ($r in $xRows) {
    $ws.Cells.Item($r, 30).Value = "X"
    $ws.Cells.Item($r, 31).Value = "X"
}

Write-Output "done"
